$wb = $excel.ActiveWorkbook

# Use "SearchInputFields" as the style donor - it already has the exact
# header/legend styling (2,40,41,42) that the new SearchAlias sheet needs.
$src = $wb.Worksheets.Item("SearchInputFields")

# New sheet is inserted right before "WorkBasketResultFields" (i.e. right
# after "SearchResultFields"), matching the commit's sheet ordering.
$target = $wb.Worksheets.Item("WorkBasketResultFields")
$new = $wb.Worksheets.Add($target)
$new.Name = "SearchAlias"

# Copy cell formatting (styles) from the donor sheet.
$src.Range("A1:D1").Copy()
$new.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$src.Range("A2:F3").Copy()
$new.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Row 3's A3/B3 need the same style as C3/D3 (44) rather than the donor's
# default (10), so re-paste that specific format onto them.
$src.Range("C3").Copy()
$new.Range("A3").PasteSpecial(-4122)
$new.Range("B3").PasteSpecial(-4122)

# Populate the cell values.
$new.Range("A1").Value = "SearchAlias"
$new.Range("B1").Value = "PrimaryKeyInRed"
$new.Range("C1").Value = "PrimaryAndForeignKey Orange"
$new.Range("D1").Value = "ForeignKey Brown"
$new.Range("B2").Value = "Unique alias id for a case field of a case type"
$new.Range("C2").Value = "This is just the CaseField ID for top level fields, or object notation pointing to a complex type field e.g. applicantAddress.AddressLine1"
$new.Range("A3").Value = "CaseTypeID"
$new.Range("B3").Value = "SearchAliasID"
$new.Range("C3").Value = "CaseFieldID"

# Row heights (matches donor, set explicitly to be safe).
$new.Rows.Item(1).RowHeight = 18
$new.Rows.Item(2).RowHeight = 113

# Column widths.
$new.Columns.Item(1).ColumnWidth = 13.7501
$new.Columns.Item(2).ColumnWidth = 14.9167
$new.Columns.Item(3).ColumnWidth = 17.5834

# Select B2 on the new sheet, making it the active sheet/cell - this also
# naturally clears tabSelected on whichever sheet was active before.
$new.Range("B2").Select()
